$d = $word.ActiveDocument
$word.UserName = "USER"

function Find-Text($range, $text) {
    $f = $range.Find
    [void]$f.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $f.Found) {
        throw "Could not find text: $text"
    }
    return $range
}

# ---------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark near the top of the document
#    (<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
#    right before "Interactivo F1: ...").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "CN_06_03_REC150_IMG01" -> "CN_06_03_REC140_IMG01"
#    Only the single digit "5" -> "4" is a tracked change, the common
#    prefix "CN_06_03_REC1" and suffix "0_IMG01" stay untouched.
# ---------------------------------------------------------------------
$r1 = $d.Content
Find-Text $r1 "CN_06_03_REC150_IMG01" | Out-Null
$start1 = $r1.Start
$digit1 = $d.Range($start1 + 13, $start1 + 14)
if ($digit1.Text -ne "5") { throw "Unexpected text at digit1: $($digit1.Text)" }
$digit1.Text = "4"

# ---------------------------------------------------------------------
# 3) "CN_06_03_REC150" (immediately followed by a separate "_IMG02" run)
#    -> "CN_06_03_REC140", same single-digit tracked change, and the
#    document's _GoBack bookmark is re-dropped right after this edit
#    (between the inserted "4" and the deleted "5").
# ---------------------------------------------------------------------
$r2 = $d.Content
Find-Text $r2 "CN_06_03_REC150_IMG02" | Out-Null
$start2 = $r2.Start
$digit2 = $d.Range($start2 + 13, $start2 + 14)
if ($digit2.Text -ne "5") { throw "Unexpected text at digit2: $($digit2.Text)" }
$digit2.Text = "4"
$bm2 = $d.Range($start2 + 14, $start2 + 14)
$d.Bookmarks.Add("_GoBack", $bm2)

# ---------------------------------------------------------------------
# 4) "CN_06_03_REC150" (immediately followed by a separate "_IMG03" run)
#    -> "CN_06_03_REC140", same single-digit tracked change.
# ---------------------------------------------------------------------
$r3 = $d.Content
Find-Text $r3 "CN_06_03_REC150_IMG03" | Out-Null
$start3 = $r3.Start
$digit3 = $d.Range($start3 + 13, $start3 + 14)
if ($digit3.Text -ne "5") { throw "Unexpected text at digit3: $($digit3.Text)" }
$digit3.Text = "4"

Write-Host "done"
